$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.140.32"
$c.ClearFormats()
$ws.Range("E2").Value = "  -3.35%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.861.77"
$c.ClearFormats()
$ws.Range("E3").Value = "  -4.12%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.90"
$c.ClearFormats()
$ws.Range("E5").Value = "  -3.39%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4659"
$c.ClearFormats()
$ws.Range("E7").Value = "  -2.88%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2823"
$c.ClearFormats()
$ws.Range("E8").Value = "  -3.07%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06551"
$c.ClearFormats()
$ws.Range("E9").Value = "  -3.50%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.09"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.77%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07812"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.48%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "96.22"
$c.ClearFormats()
$ws.Range("E12").Value = "  -7.64%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.857.56"
$c.ClearFormats()
$ws.Range("E13").Value = "  -4.40%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.125"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.32%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6692"
$c.ClearFormats()
$ws.Range("E15").Value = "  -3.87%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "280.53"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.34%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.157.18"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.25%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.09%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.473"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.63%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.61"
$c.ClearFormats()
$ws.Range("E20").Value = "  -2.94%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.098.63"
$c.ClearFormats()
$ws.Range("E21").Value = "  -5.00%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.000007243"
$c.ClearFormats()
$ws.Range("E22").Value = "  -4.84%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.9992"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.136"
$c.ClearFormats()
$ws.Range("E24").Value = "  -4.76%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.313"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.64%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "164.97"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.47%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.87"
$c.ClearFormats()
$ws.Range("E27").Value = "  -5.00%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.901"
$c.ClearFormats()
$ws.Range("E28").Value = "  -9.87%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.342"
$c.ClearFormats()
$ws.Range("E29").Value = "  -4.34%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09594"
$c.ClearFormats()
$ws.Range("E30").Value = "  -4.86%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.393"
$c.ClearFormats()
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("E32").Value = "  -4.39%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.098"
$c.ClearFormats()
$ws.Range("E33").Value = "  -5.85%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04651"
$c.ClearFormats()
$ws.Range("E34").Value = "  -3.76%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.6998"
$c.ClearFormats()
$ws.Range("E35").Value = "  -5.58%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.094"
$c.ClearFormats()
$ws.Range("E36").Value = "  -3.21%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.709"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  -5.47%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.294"
$c.ClearFormats()
$ws.Range("E39").Value = "  -6.93%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.521"
$c.ClearFormats()
$ws.Range("E40").Value = "  -4.76%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "72.73"
$c.ClearFormats()
$ws.Range("E41").Value = "  -5.21%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.8521"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.12%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.920"
$c.ClearFormats()
$ws.Range("E43").Value = "  -5.46%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.08%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4158"
$c.ClearFormats()
$ws.Range("E45").Value = "  -5.00%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "103.39"
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "990.10"
$c.ClearFormats()
$ws.Range("E47").Value = "  -3.51%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.168"
$c.ClearFormats()
$ws.Range("E48").Value = "  -5.37%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.189"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.33%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "34.12"
$c.ClearFormats()
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("E51").Value = "  -6.02%  "
